$wb = $excel.ActiveWorkbook

$wsInstrucciones = $wb.Worksheets.Item(1)
$wsVelocidad     = $wb.Worksheets.Item(2)
$wsEstimacion    = $wb.Worksheets.Item(3)

# --- Sheet "VelocidadEquipo" ---
$wsVelocidad.Range("A4").Value = "N.º Horas por Semana"
$wsVelocidad.Range("A5").Value = "N.º Semanas"
$wsVelocidad.Range("A6").Value = "N.º de Desarrolladores"

# --- Sheet "Instrucciones" ---
$c12Text = "Consignen la lista de actividades transversales que se deben realizar durante el sprint y los tiempos estimados para las mismas. Una tarea transversal es una tarea que se debe realizar y no se hace para una historia particular. Un ejemplo de tarea transversal es la actividad de grooming de historias."
$wsInstrucciones.Range("C12").Value = $c12Text

$wsInstrucciones.Range("B9").Value = "N.º Horas por semana"

# --- Sheet "EstimaciónHUBase" ---
$wsEstimacion.Range("A3").Value = "Título"
$wsEstimacion.Range("A6").Value = "Sumatoria de horas"
$wsEstimacion.Range("C9").Value = "Estimación de horas"

# --- back to "Instrucciones" ---
$c9Text = "Horas por semana dedicadas al proyecto por cada miembro. Tener en cuenta que:" + [char]10 + `
    "- La dedicación esperada por semana para el curso según el número de créditos es de 12 horas." + [char]10 + `
    "- Se espera que los estudiantes dediquen entre 3 y 4 horas a actividades de aprendizaje que no hacen parte del proyecto."
$wsInstrucciones.Range("C9").Value = $c9Text

$wsInstrucciones.Range("B10").Value = "N.º Semanas"
$wsInstrucciones.Range("B11").Value = "N.º de Desarrolladores"

# --- Nudge the header images very slightly (matches the tiny reposition/resize
#     that Excel recorded for these two pictures) ---
$picInstrucciones = $wsInstrucciones.Shapes.Item(1)
$picInstrucciones.Left = 1.5
$picInstrucciones.Width = 594.0
$picInstrucciones.Height = 139.2187401574803

$picVelocidad = $wsVelocidad.Shapes.Item(1)
$picVelocidad.Left = 0.75
$picVelocidad.Width = 513.0
$picVelocidad.Height = 120.23433070866142

# --- Restore the selections shown in each sheet when the workbook was last saved ---
$wsVelocidad.Activate()
$wsVelocidad.Range("G4").Select()

$wsEstimacion.Activate()
$wsEstimacion.Range("E9").Select()

$wsInstrucciones.Activate()
$wsInstrucciones.Range("E9").Select()
